$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3743
$ws.Range("J3").Value = 3932
$ws.Range("J4").Value = 875
$ws.Range("I5").Value = 718
$ws.Range("J5").Value = 314
$ws.Range("J6").Value = 4616
$ws.Range("I7").Value = 26212
$ws.Range("J7").Value = 13480

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 39
$ws.Range("J7").Value = 399
$ws.Range("J8").Value = 870
$ws.Range("J11").Value = 199
$ws.Range("J15").Value = 153
$ws.Range("J16").Value = 39
$ws.Range("J19").Value = 400
$ws.Range("J20").Value = 271
$ws.Range("J25").Value = 71
$ws.Range("J26").Value = 25
$ws.Range("J29").Value = 772
$ws.Range("J30").Value = 56
$ws.Range("J31").Value = 109
$ws.Range("J33").Value = 613
$ws.Range("J36").Value = 198
$ws.Range("J37").Value = 423
$ws.Range("J40").Value = 25
$ws.Range("J41").Value = 83
$ws.Range("J42").Value = 519
$ws.Range("J46").Value = 49
$ws.Range("J52").Value = 370
$ws.Range("J53").Value = 127
$ws.Range("I63").Value = 224
$ws.Range("J64").Value = 90
$ws.Range("J65").Value = 355
$ws.Range("J66").Value = 40
$ws.Range("J67").Value = 519
$ws.Range("J72").Value = 55
$ws.Range("J73").Value = 119
$ws.Range("J76").Value = 194
$ws.Range("J78").Value = 182
$ws.Range("J79").Value = 390
$ws.Range("J80").Value = 24
$ws.Range("J84").Value = 121
$ws.Range("J85").Value = 610
$ws.Range("J90").Value = 157
$ws.Range("J91").Value = 157
$ws.Range("J94").Value = 119
$ws.Range("J95").Value = 211
$ws.Range("J96").Value = 150
$ws.Range("J97").Value = 91
$ws.Range("J99").Value = 195
$ws.Range("I101").Value = 26212
$ws.Range("J101").Value = 13480

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 75
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 197
$ws.Range("J5").Value = 28
$ws.Range("J6").Value = 199
$ws.Range("J7").Value = 613

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 242
$ws.Range("J3").Value = 267
$ws.Range("J4").Value = 43
$ws.Range("J5").Value = 35
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 772

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 154
$ws.Range("J7").Value = 610

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 66
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 211

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 115
$ws.Range("J7").Value = 519

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 148
$ws.Range("J6").Value = 125
$ws.Range("J7").Value = 423

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 75
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 51
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 117
$ws.Range("J6").Value = 146
$ws.Range("J7").Value = 519

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 100
$ws.Range("J6").Value = 148
$ws.Range("J7").Value = 400

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 83
$ws.Range("J7").Value = 370

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 255
$ws.Range("J3").Value = 275
$ws.Range("J4").Value = 47
$ws.Range("J6").Value = 267
$ws.Range("J7").Value = 870

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 99
$ws.Range("J3").Value = 109
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 147
$ws.Range("J7").Value = 390

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 132
$ws.Range("J4").Value = 11
$ws.Range("J6").Value = 136
$ws.Range("J7").Value = 399

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J2").Value = 6
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J3").Value = 14
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 39
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 62
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 39
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J4").Value = 5
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 25
